# Resort the sheet order: "总计" should come before "2022-Q1".
# (Data/content of each named sheet is unchanged -- only the tab order,
#  and therefore which physical sheet part each name/rId points at, moves.)

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# Move "总计" to sit immediately before "2022-Q1", giving the final order:
#   1) 总计
#   2) 2022-Q1
$totalSheet.Move($q1Sheet)
